# Apply "Now handles trailing rows" edit
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Rename the table titles so they are distinguishable (old sheet vs new sheet)
$ws1.Range("A1").Value = "Example TableOld"
$ws2.Range("A1").Value = "Example TableNew"

# Sheet1 gains a trailing row with a note, using the default/general style
$ws1.Range("A10").Value = "Trailing row here"

# Move the active selection/tab from Sheet1 to Sheet2
$ws1.Range("A1").Select()
$ws2.Activate()
$ws2.Range("A2").Select()
